$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Structural changes -----------------------------------------------
# 1) Insert a new row for "Latte" right after "Vanilla Cream" (becomes row 3).
$ws.Rows.Item(3).Insert()

# 2) Insert a new row for "Aloo Tikki" right before "Peanut Masala"
#    (Peanut Masala is currently at row 13 after the previous insert, so
#    inserting at row 13 pushes it down to 14 and creates space at 13).
$ws.Rows.Item(13).Insert()

# 3) Remove the trailing rows that no longer exist in the updated sheet:
#    Veggie Delite, Paneer Tikka, Corn n Peas, Peri Peri Chicken, Ham Chicken.
#    After the two inserts above, these currently sit at rows 20-24.
$ws.Range("20:24").Delete()

# --- Force Price/Tax columns to be stored as text (data rows only, the
#     header row in row 1 keeps its original formatting) ------------------
$ws.Range("C2:D19").NumberFormat = "@"

# --- Write the final data ------------------------------------------------
$ws.Cells.Item(2,1).Value = "Vanilla Cream"
$ws.Cells.Item(2,2).Value = 50
$ws.Cells.Item(2,3).Value = "8750.00"
$ws.Cells.Item(2,4).Value = "700.0000"

$ws.Cells.Item(3,1).Value = "Latte"
$ws.Cells.Item(3,2).Value = 20
$ws.Cells.Item(3,3).Value = "2000.00"
$ws.Cells.Item(3,4).Value = "160.0000"

$ws.Cells.Item(4,1).Value = "Filter Coffee"
$ws.Cells.Item(4,2).Value = 30
$ws.Cells.Item(4,3).Value = "2700.00"
$ws.Cells.Item(4,4).Value = "216.0000"

$ws.Cells.Item(5,1).Value = "Cold Coffee"
$ws.Cells.Item(5,2).Value = 12
$ws.Cells.Item(5,3).Value = "1800.00"
$ws.Cells.Item(5,4).Value = "144.0000"

$ws.Cells.Item(6,1).Value = "Black Coffee"
$ws.Cells.Item(6,2).Value = 10
$ws.Cells.Item(6,3).Value = "880.00"
$ws.Cells.Item(6,4).Value = "70.4000"

$ws.Cells.Item(7,1).Value = "Normal Tea"
$ws.Cells.Item(7,2).Value = 10
$ws.Cells.Item(7,3).Value = "250.00"
$ws.Cells.Item(7,4).Value = "20.0000"

$ws.Cells.Item(8,1).Value = "Lemon Tea"
$ws.Cells.Item(8,2).Value = 10
$ws.Cells.Item(8,3).Value = "350.00"
$ws.Cells.Item(8,4).Value = "28.0000"

$ws.Cells.Item(9,1).Value = "Ginger Tea"
$ws.Cells.Item(9,2).Value = 10
$ws.Cells.Item(9,3).Value = "400.00"
$ws.Cells.Item(9,4).Value = "32.0000"

$ws.Cells.Item(10,1).Value = "Honey Tea"
$ws.Cells.Item(10,2).Value = 10
$ws.Cells.Item(10,3).Value = "405.00"
$ws.Cells.Item(10,4).Value = "32.4000"

$ws.Cells.Item(11,1).Value = "Bhel Puri"
$ws.Cells.Item(11,2).Value = 10
$ws.Cells.Item(11,3).Value = "2200.00"
$ws.Cells.Item(11,4).Value = "264.0000"

$ws.Cells.Item(12,1).Value = "Sev Puri"
$ws.Cells.Item(12,2).Value = 10
$ws.Cells.Item(12,3).Value = "1700.00"
$ws.Cells.Item(12,4).Value = "204.0000"

$ws.Cells.Item(13,1).Value = "Aloo Tikki"
$ws.Cells.Item(13,2).Value = 10
$ws.Cells.Item(13,3).Value = "2800.00"
$ws.Cells.Item(13,4).Value = "336.0000"

$ws.Cells.Item(14,1).Value = "Peanut Masala"
$ws.Cells.Item(14,2).Value = 10
$ws.Cells.Item(14,3).Value = "2990.00"
$ws.Cells.Item(14,4).Value = "358.8000"

$ws.Cells.Item(15,1).Value = "Chilli Cheese"
$ws.Cells.Item(15,2).Value = 10
$ws.Cells.Item(15,3).Value = "3450.00"
$ws.Cells.Item(15,4).Value = "414.0000"

$ws.Cells.Item(16,1).Value = "Chicken Tikki"
$ws.Cells.Item(16,2).Value = 10
$ws.Cells.Item(16,3).Value = "2980.00"
$ws.Cells.Item(16,4).Value = "536.4000"

$ws.Cells.Item(17,1).Value = "Kheema chat"
$ws.Cells.Item(17,2).Value = 10
$ws.Cells.Item(17,3).Value = "3100.00"
$ws.Cells.Item(17,4).Value = "558.0000"

$ws.Cells.Item(18,1).Value = "Shambi Kabab"
$ws.Cells.Item(18,2).Value = 10
$ws.Cells.Item(18,3).Value = "3750.00"
$ws.Cells.Item(18,4).Value = "675.0000"

$ws.Cells.Item(19,1).Value = "Chicken Tandoor"
$ws.Cells.Item(19,2).Value = 10
$ws.Cells.Item(19,3).Value = "3400.00"
$ws.Cells.Item(19,4).Value = "612.0000"
